$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 with the 2021 figures, matching the style of the
# existing year-label cells in column A (row 11).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A12").Value = "2021年"

$ws.Range("B12").Value = 59
$ws.Range("D12").Value = 1137
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 1076
